# Update the "Fitness" (best-fitness-so-far) values in column C for rows 2-252
# (Generation 0 through 250) to reflect the results of the re-run logged for run_18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fitnessValues = @(
    12744, 12510, 12138, 12138, 11077, 11077, 9946, 9392, 9392, 9392, 9020, 9020, 9020, 9020, 8593,
    8593, 8593, 8593, 8593, 8593, 8593, 8593, 8593, 8593, 8593, 8593, 8593, 8593, 8593, 8593,
    8593, 8593, 8593, 8593, 8593, 8593, 8593, 8593, 8593, 8053, 8053, 8053, 8053, 8053, 8053,
    8053, 8053, 8053, 8053, 8028, 8028, 8028, 8028, 8028, 8028, 8028, 8028, 8028, 8028, 8028,
    8028, 8028, 8028, 8028, 8028, 8028, 8028, 8028, 8028, 8028, 8028, 8028, 8028, 8028, 8006,
    8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006,
    8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006,
    8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006,
    8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006,
    8006, 8006, 8006, 8006, 8006, 8006, 8006, 8006, 7569, 7569, 7569, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569
)

$startRow = 2
for ($i = 0; $i -lt $fitnessValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $fitnessValues[$i]
}
